# Applies scheduled-runner market-data refresh to all 8 job sheets.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N)
# for the rows whose source data changed.
$wb = $excel.ActiveWorkbook

# ---- Sheet 1: ALC ----
$ws = $wb.Worksheets.Item(1)
$ws.Range("H53").Value = 2640.5833
$ws.Range("I53").Value = 88.583336
$ws.Range("J53").Value = 5192.5835
$ws.Range("K53").Value = 88.583336
$ws.Range("L53").Value = 5192.5835
$ws.Range("M53").Value = 548.416664
$ws.Range("N53").Value = -6466.5835
$ws.Range("H113").Value = 66669868
$ws.Range("I113").Value = 1000000000
$ws.Range("J113").Value = 3428.5715
$ws.Range("K113").Value = 1000000000
$ws.Range("L113").Value = 3428.5715
$ws.Range("M113").Value = -999996746
$ws.Range("N113").Value = -9936.5715
$ws.Range("H116").Value = 34383150
$ws.Range("J116").Value = 37041320
$ws.Range("L116").Value = 37041320
$ws.Range("N116").Value = -37048204
$ws.Range("H132").Value = 2855.5916
$ws.Range("I132").Value = 2720.5688
$ws.Range("J132").Value = 3458
$ws.Range("K132").Value = 8161.7064
$ws.Range("L132").Value = 10374
$ws.Range("M132").Value = -5631.7064
$ws.Range("N132").Value = -15434
$ws.Range("H137").Value = 2196
$ws.Range("I137").Value = 1502.8846
$ws.Range("J137").Value = 2453.4429
$ws.Range("K137").Value = 4508.6538
$ws.Range("L137").Value = 7360.3287
$ws.Range("M137").Value = -1958.6538
$ws.Range("N137").Value = -12460.3287

# ---- Sheet 2: ARM ----
$ws = $wb.Worksheets.Item(2)
$ws.Range("H2").Value = 1550
$ws.Range("J2").Value = 2004.5
$ws.Range("L2").Value = 2004.5
$ws.Range("N2").Value = -2230.5
$ws.Range("H32").Value = 110454.664
$ws.Range("I32").Value = 118663.984
$ws.Range("K32").Value = 118663.984
$ws.Range("M32").Value = -118376.984
$ws.Range("H45").Value = 3244.4167
$ws.Range("I45").Value = 3102.7144
$ws.Range("J45").Value = 3442.8
$ws.Range("K45").Value = 3102.7144
$ws.Range("L45").Value = 3442.8
$ws.Range("M45").Value = -2725.7144
$ws.Range("N45").Value = -4196.8
$ws.Range("H74").Value = 3185
$ws.Range("I74").Value = 3665.889
$ws.Range("K74").Value = 3665.889
$ws.Range("M74").Value = -2791.889
$ws.Range("H77").Value = 3185
$ws.Range("I77").Value = 3665.889
$ws.Range("K77").Value = 18329.445
$ws.Range("M77").Value = -13961.445
$ws.Range("H97").Value = 261.57895
$ws.Range("I97").Value = 225.5625
$ws.Range("J97").Value = 453.66666
$ws.Range("K97").Value = 225.5625
$ws.Range("L97").Value = 453.66666
$ws.Range("M97").Value = 270.4375
$ws.Range("N97").Value = -1445.66666
$ws.Range("H110").Value = 66668064
$ws.Range("I110").Value = 71429930
$ws.Range("J110").Value = 2000
$ws.Range("K110").Value = 71429930
$ws.Range("L110").Value = 2000
$ws.Range("M110").Value = -71427885
$ws.Range("N110").Value = -6090
$ws.Range("H116").Value = 1550
$ws.Range("J116").Value = 2004.5
$ws.Range("L116").Value = 2004.5
$ws.Range("N116").Value = -6592.5
$ws.Range("H122").Value = 3216.7878
$ws.Range("I122").Value = 1637.0476
$ws.Range("K122").Value = 4911.142800000001
$ws.Range("M122").Value = -2461.142800000001

# ---- Sheet 3: BSM ----
$ws = $wb.Worksheets.Item(3)
$ws.Range("H3").Value = 1550
$ws.Range("J3").Value = 2004.5
$ws.Range("L3").Value = 2004.5
$ws.Range("N3").Value = -2232.5
$ws.Range("H105").Value = 1910.6923
$ws.Range("I105").Value = 1751.9333
$ws.Range("K105").Value = 1751.9333
$ws.Range("M105").Value = -4.933299999999917
$ws.Range("I107").Value = 8242.931
$ws.Range("K107").Value = 8242.931
$ws.Range("M107").Value = -6322.931
$ws.Range("H141").Value = 72499.5
$ws.Range("J141").Value = 72499.5
$ws.Range("L141").Value = 72499.5
$ws.Range("N141").Value = -82859.5

# ---- Sheet 4: CRP ----
$ws = $wb.Worksheets.Item(4)
$ws.Range("H16").Value = 3367
$ws.Range("I16").Value = 3276.2778
$ws.Range("K16").Value = 3276.2778
$ws.Range("M16").Value = -2989.2778
$ws.Range("H31").Value = 2660.7273
$ws.Range("I31").Value = 838.7273
$ws.Range("K31").Value = 838.7273
$ws.Range("M31").Value = -543.7273
$ws.Range("H34").Value = 2660.7273
$ws.Range("I34").Value = 838.7273
$ws.Range("K34").Value = 838.7273
$ws.Range("M34").Value = -636.7273
$ws.Range("H113").Value = 3367
$ws.Range("I113").Value = 3276.2778
$ws.Range("K113").Value = 3276.2778
$ws.Range("M113").Value = -1106.2778
$ws.Range("H122").Value = 2692
$ws.Range("I122").Value = 1786.4166
$ws.Range("K122").Value = 5359.2498
$ws.Range("M122").Value = -2909.2498
$ws.Range("H132").Value = 2109.5862
$ws.Range("I132").Value = 1969.5555
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 5908.666499999999
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -3378.666499999999
$ws.Range("N132").Value = -17060

# ---- Sheet 5: CUL ----
$ws = $wb.Worksheets.Item(5)
$ws.Range("H34").Value = 659.3570999999999
$ws.Range("I34").Value = 147.28572
$ws.Range("J34").Value = 1171.4286
$ws.Range("K34").Value = 441.85716
$ws.Range("L34").Value = 3514.2858
$ws.Range("M34").Value = -357.85716
$ws.Range("N34").Value = -3682.2858
$ws.Range("H113").Value = 581.8
$ws.Range("I113").Value = 341
$ws.Range("J113").Value = 669.36365
$ws.Range("K113").Value = 1023
$ws.Range("L113").Value = 2008.09095
$ws.Range("M113").Value = 1147
$ws.Range("N113").Value = -6348.09095
$ws.Range("H127").Value = 999.3333
$ws.Range("J127").Value = 999.3333
$ws.Range("L127").Value = 2997.9999
$ws.Range("N127").Value = -12917.9999
$ws.Range("H129").Value = 2094.0908
$ws.Range("I129").Value = 1482.5
$ws.Range("J129").Value = 2323.4375
$ws.Range("K129").Value = 4447.5
$ws.Range("L129").Value = 6970.3125
$ws.Range("M129").Value = 552.5
$ws.Range("N129").Value = -16970.3125
$ws.Range("H131").Value = 10324.24
$ws.Range("J131").Value = 13109.421
$ws.Range("L131").Value = 39328.263
$ws.Range("N131").Value = -49408.263
$ws.Range("H140").Value = 2643
$ws.Range("I140").Value = 2643
$ws.Range("K140").Value = 7929
$ws.Range("M140").Value = -2749

# ---- Sheet 6: GSM ----
$ws = $wb.Worksheets.Item(6)
$ws.Range("H80").Value = 2206.3
$ws.Range("I80").Value = 2184.2856
$ws.Range("J80").Value = 2213
$ws.Range("K80").Value = 2184.2856
$ws.Range("L80").Value = 2213
$ws.Range("M80").Value = -1186.2856
$ws.Range("N80").Value = -4209
$ws.Range("H83").Value = 2206.3
$ws.Range("I83").Value = 2184.2856
$ws.Range("J83").Value = 2213
$ws.Range("K83").Value = 10921.428
$ws.Range("L83").Value = 11065
$ws.Range("M83").Value = -5929.428
$ws.Range("N83").Value = -21049
$ws.Range("H102").Value = 2572.261
$ws.Range("I102").Value = 1609
$ws.Range("K102").Value = 1609
$ws.Range("M102").Value = 13
$ws.Range("H113").Value = 4500
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 4500
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 4500
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -8840
$ws.Range("H122").Value = 3236.0667
$ws.Range("I122").Value = 2961.75
$ws.Range("J122").Value = 4333.3335
$ws.Range("K122").Value = 8885.25
$ws.Range("L122").Value = 13000.0005
$ws.Range("M122").Value = -6435.25
$ws.Range("N122").Value = -17900.0005
$ws.Range("H126").Value = 9128.429
$ws.Range("I126").Value = 13379
$ws.Range("K126").Value = 40137
$ws.Range("M126").Value = -37667

# ---- Sheet 7: LTW ----
$ws = $wb.Worksheets.Item(7)
$ws.Range("H7").Value = 55559000
$ws.Range("I7").Value = 500000000
$ws.Range("K7").Value = 500000000
$ws.Range("M7").Value = -499999888
$ws.Range("H40").Value = 3950.5
$ws.Range("I40").Value = 3000
$ws.Range("J40").Value = 4584.1665
$ws.Range("K40").Value = 3000
$ws.Range("L40").Value = 4584.1665
$ws.Range("M40").Value = -2864
$ws.Range("N40").Value = -4856.1665
$ws.Range("H46").Value = 2641.5
$ws.Range("J46").Value = 3033.3333
$ws.Range("L46").Value = 3033.3333
$ws.Range("N46").Value = -3409.3333
$ws.Range("H93").Value = 1025.1708
$ws.Range("I93").Value = 1044
$ws.Range("J93").Value = 947.5
$ws.Range("K93").Value = 1044
$ws.Range("L93").Value = 947.5
$ws.Range("M93").Value = 204
$ws.Range("N93").Value = -3443.5
$ws.Range("H122").Value = 4120
$ws.Range("I122").Value = 3750
$ws.Range("J122").Value = 4490
$ws.Range("K122").Value = 11250
$ws.Range("L122").Value = 13470
$ws.Range("M122").Value = -8800
$ws.Range("N122").Value = -18370
$ws.Range("H126").Value = 55559000
$ws.Range("I126").Value = 500000000
$ws.Range("K126").Value = 1500000000
$ws.Range("M126").Value = -1499997530
$ws.Range("H132").Value = 6618.364
$ws.Range("I132").Value = 3707.6155
$ws.Range("J132").Value = 10822.777
$ws.Range("K132").Value = 11122.8465
$ws.Range("L132").Value = 32468.331
$ws.Range("M132").Value = -8592.8465
$ws.Range("N132").Value = -37528.331

# ---- Sheet 8: WVR ----
$ws = $wb.Worksheets.Item(8)
$ws.Range("H54").Value = 10000
$ws.Range("J54").Value = 10000
$ws.Range("L54").Value = 10000
$ws.Range("N54").Value = -11040
$ws.Range("H96").Value = 4562.6665
$ws.Range("I96").Value = 4231.625
$ws.Range("J96").Value = 5224.75
$ws.Range("K96").Value = 4231.625
$ws.Range("L96").Value = 5224.75
$ws.Range("M96").Value = -2858.625
$ws.Range("N96").Value = -7970.75
$ws.Range("H107").Value = 377.15
$ws.Range("I107").Value = 370.5263
$ws.Range("K107").Value = 1111.5789
$ws.Range("M107").Value = 808.4211
$ws.Range("H113").Value = 11944
$ws.Range("I113").Value = 3888.5
$ws.Range("K113").Value = 11665.5
$ws.Range("M113").Value = -9495.5
$ws.Range("H132").Value = 479884.25
$ws.Range("I132").Value = 693050.7
$ws.Range("K132").Value = 2079152.1
$ws.Range("M132").Value = -2076622.1
$ws.Range("H136").Value = 7467.8945
$ws.Range("I136").Value = 6126.4
$ws.Range("K136").Value = 18379.2
$ws.Range("M136").Value = -15829.2
